$d = $word.ActiveDocument

function Get-ParaIndexAt($pos) {
  $idx = 0
  $result = 0
  foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {
      $result = $idx
    }
  }
  return $result
}

# 1) Collapse the split-run "Sensor <Item> – " headings (each was typed as
#    several separate runs: "Sensor " / "<Item>" / " – ") into a single bold
#    run, and likewise collapse the two-run sentence that follows it into a
#    single normal run. A same-text Find/Replace across the run boundaries
#    rebuilds the matched span as freshly-merged runs.
$d.Content.Find.Execute("Sensor Cacerola", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sensor Cacerola", 2) | Out-Null
$d.Content.Find.Execute("Sensor Wok", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sensor Wok", 2) | Out-Null
$d.Content.Find.Execute("Sensor Paellera", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sensor Paellera", 2) | Out-Null

# 2) The "PLANTILLA MANIOBRA" heading was split as "PLANT" + a stray
#    "_GoBack" bookmark + "ILLA MANIOBRA". Re-merge into one run; the
#    Find/Replace also drops the now-redundant bookmark from here (it gets
#    re-added at the very end of the document below, mirroring Word's own
#    habit of relocating "_GoBack" to the most recent edit point).
$d.Content.Find.Execute("PLANTILLA MANIOBRA", $true, $false, $false, $false, $false,
                         $true, 1, $false, "PLANTILLA MANIOBRA", 2) | Out-Null

# 3) Remove the entire "PLANTILLA START/STOP" section: from the page-break
#    paragraph that introduces it through to its last content paragraph,
#    leaving only the trailing blank paragraph that precedes the sectPr.
$headingRange = $d.Content
$headingRange.Find.Execute("PLANTILLA START/STOP", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0) | Out-Null
$headingIdx = Get-ParaIndexAt($headingRange.Start)
$pageBreakIdx = $headingIdx - 1
$startPos = $d.Paragraphs.Item($pageBreakIdx).Range.Start

$endRange = $d.Content
$endRange.Find.Execute("realiza cuatro.", $true, $false, $false, $false, $false,
                        $true, 1, $false, "", 0) | Out-Null
$endParaIdx = Get-ParaIndexAt($endRange.Start)
$endPos = $d.Paragraphs.Item($endParaIdx).Range.End

$delRange = $d.Range($startPos, $endPos)
$delRange.Delete()

# The single remaining paragraph loses its justification and becomes the
# new home of the "_GoBack" bookmark (an empty paragraph holding only the
# bookmark, right before the section properties).
$last = $d.Paragraphs.Last
$last.Alignment = 0
$d.Bookmarks.Add("_GoBack", $last.Range) | Out-Null
